# edit.ps1 - apply the OOXML changes described by the target diff:
#   1. Slide 3 ("Functionality" slide): swap the first two bullet lines so that
#      "add a new customer(s)" now comes before "search customers by name"
#      (the first bullet becomes the plural "add a new customers").
#   2. Reorder slides 4 and 5 so the "add a new customer" title slide now
#      appears before the "search customers by name" title slide.
#   3. (best effort) Register an empty presentation-level slide-guide list,
#      matching the tiny <p:extLst>/<p15:sldGuideLst/> stub added by
#      PowerPoint when the Guides feature is touched.

$p = $ppt.ActivePresentation

# --- 1. Slide 3: swap the two bullet paragraphs in the content placeholder ---
$s3 = $p.Slides.Item(3)
$contentShape = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.Name -eq "Content Placeholder 2") {
        $contentShape = $sh
    }
}
if ($contentShape -eq $null) {
    $contentShape = $s3.Shapes.Item(5)
}
$tr = $contentShape.TextFrame.TextRange

# Original paragraph 1: "search customers by name"
# Original paragraph 2: "add a new customer"
# Update the later paragraph first so the earlier paragraph's range/offset
# is not disturbed by the text-length change.
$tr.Paragraphs(2, 1).Text = "search customers by name"
$tr.Paragraphs(1, 1).Text = "add a new customers"

# --- 2. Swap the display order of the two title slides (positions 4 and 5) ---
$p.Slides.Item(5).MoveTo(4)

# --- 3. Best-effort: touch the presentation-level Guides collection so an
#        (empty) slide guide list gets registered, matching the stub extLst
#        PowerPoint writes out after the Guides feature is used. ---
try {
    $guides = $p.Guides
    $guides.Add(1, 3000) | Out-Null
} catch {
}
